# Apply the price / volume updates from the Sept 7 2024 GitHub Actions crypto refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell([string]$addr, [string]$text) {
    # Leading apostrophe forces Excel to store numeric-looking strings (e.g. '1.00')
    # as literal text instead of silently coercing them to numbers.
    $ws.Range($addr).Value = "'" + $text
    # Resetting the style strips the 'number stored as text' quote-prefix flag Excel
    # would otherwise attach, keeping formatting identical to the original cell.
    $ws.Range($addr).Style = "Normal"
}

Set-TextCell "D2" "54.264.76"
Set-TextCell "E2" "  -2.74%  "
Set-TextCell "D3" "2.258.99"
Set-TextCell "E3" "  -3.70%  "
Set-TextCell "D4" "0.999"
Set-TextCell "E4" "  -0.17%  "
Set-TextCell "D5" "494.18"
Set-TextCell "E5" "  -2.04%  "
Set-TextCell "D6" "127.30"
Set-TextCell "E6" "  -1.48%  "
Set-TextCell "E7" "  +0.10%  "
Set-TextCell "E8" "  -1.75%  "
Set-TextCell "D9" "2.283.84"
Set-TextCell "E9" "  -3.05%  "
Set-TextCell "D10" "0.0940"
Set-TextCell "E10" "  -3.33%  "
Set-TextCell "E11" "  +0.31%  "
Set-TextCell "E12" "  +0.59%  "
Set-TextCell "D13" "4.63"
Set-TextCell "E13" "  -3.10%  "
Set-TextCell "D14" "2.686.90"
Set-TextCell "E14" "  -2.74%  "
Set-TextCell "D15" "21.50"
Set-TextCell "E15" "  -0.47%  "
Set-TextCell "D16" "54.127.24"
Set-TextCell "E16" "  -2.93%  "
Set-TextCell "E17" "  -2.18%  "
Set-TextCell "D18" "2.290.51"
Set-TextCell "E18" "  -2.43%  "
Set-TextCell "D19" "9.90"
Set-TextCell "E19" "  -0.25%  "
Set-TextCell "E20" "  +1.10%  "
Set-TextCell "D21" "301.78"
Set-TextCell "E21" "  -2.93%  "
Set-TextCell "E22" "  +3.00%  "
Set-TextCell "D23" "1.00"
Set-TextCell "E23" "  +0.12%  "
Set-TextCell "D24" "5.38"
Set-TextCell "E24" "  -1.71%  "
Set-TextCell "D25" "63.71"
Set-TextCell "E25" "  -2.30%  "
Set-TextCell "E26" "  +0.23%  "
Set-TextCell "D27" "0.372"
Set-TextCell "E27" "  +0.53%  "
Set-TextCell "D28" "2.386.86"
Set-TextCell "E28" "  -2.70%  "
Set-TextCell "E29" "  +1.81%  "
Set-TextCell "D30" "7.11"
Set-TextCell "E30" "  +0.27%  "
Set-TextCell "D31" "165.30"
Set-TextCell "E31" "  -3.52%  "
Set-TextCell "E32" "  -2.80%  "
Set-TextCell "D33" "0.0₃0681"
Set-TextCell "E33" "  -3.21%  "
Set-TextCell "E34" "  +1.90%  "
Set-TextCell "D36" "0.995"
Set-TextCell "E36" "  -0.19%  "
Set-TextCell "E37" "  +0.48%  "
Set-TextCell "D38" "17.56"
Set-TextCell "E38" "  -0.71%  "
Set-TextCell "E39" "  +1.19%  "
Set-TextCell "D40" "0.870"
Set-TextCell "E40" "  +4.77%  "
Set-TextCell "E41" "  -0.43%  "
Set-TextCell "D42" "35.42"
Set-TextCell "E42" "  -1.64%  "
Set-TextCell "D43" "0.373"
Set-TextCell "E43" "  +0.58%  "
Set-TextCell "E44" "  +0.65%  "
Set-TextCell "E45" "  -0.13%  "
Set-TextCell "D46" "126.01"
Set-TextCell "E46" "  -0.26%  "
Set-TextCell "D47" "4.79"
Set-TextCell "E47" "  -2.07%  "
Set-TextCell "D48" "0.0887"
Set-TextCell "E48" "  -0.50%  "
Set-TextCell "D49" "0.544"
Set-TextCell "E49" "  -2.12%  "
Set-TextCell "D50" "236.90"
Set-TextCell "E50" "  -0.77%  "
Set-TextCell "E51" "  +0.60%  "
